# Insert a new data row at row 85 (pushing existing rows 85-174 down to 86-175)
# and populate it with the new weekly price record for Jengibre.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(85).Insert()

$ws.Range("A85").Value = 8
$ws.Range("B85").Value = "Terminal La Palmera de La Serena"
$ws.Range("C85").Value = "Coquimbo"
$ws.Range("D85").Value = 45159
$ws.Range("E85").Value = 4
$ws.Range("F85").Value = 100114007
$ws.Range("G85").Value = "Jengibre"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 360
$ws.Range("K85").Value = 18000
$ws.Range("L85").Value = 19000
$ws.Range("M85").Value = 18500
$ws.Range("N85").Value = "$/caja 13 kilos"
$ws.Range("O85").Value = "Perú"
$ws.Range("P85").Value = 1423
$ws.Range("Q85").Value = 13
$ws.Range("R85").Value = "Hortaliza"
